# B6-PowerPoint.pptx edit (Wed, Jun 03, 2020  2:07:15 PM)
#
# 1. Three tables (on the slides that hold the "Table_0" style) are
#    switched from the custom table style
#    {73A27EE6-BACD-4465-966C-D2473F10C582} to the built-in
#    "No Style, Table Grid" style {F14ADAA1-DF4D-440E-A0C8-5A853C2F67F7}.
# 2. The deck's theme colour palette is swapped: the palette that is
#    actually applied to the slide master / all slides (the "Red Violet"
#    / Integral scheme) is replaced with the plain "Office" colour
#    scheme, mirroring the author's swap of the two theme parts.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table still using the old custom GUID ----------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq "{73A27EE6-BACD-4465-966C-D2473F10C582}") {
                $tbl.ApplyStyle("{F14ADAA1-DF4D-440E-A0C8-5A853C2F67F7}")
            }
        }
    }
}

# --- 2. Swap the theme colour scheme -----------------------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (MsoThemeColorSchemeIndex order)
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($k = 1; $k -le 12; $k++) {
    $tcs.Colors($k).RGB = $officeColors[$k - 1]
}
